$wb = $excel.ActiveWorkbook

# Sheet ALC, row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1702.0555
$ws.Range("J28").Value = 4299
$ws.Range("L28").Value = 4299
$ws.Range("N28").Value = -5269

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 233336290
$ws.Range("J51").Value = 100004450
$ws.Range("L51").Value = 100004450
$ws.Range("N51").Value = -100005418

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1994.8572
$ws.Range("I106").Value = 1878.5264
$ws.Range("K106").Value = 1878.5264
$ws.Range("M106").Value = -1247.5264

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5933.273
$ws.Range("I137").Value = 1568.0714
$ws.Range("J137").Value = 13572.375
$ws.Range("K137").Value = 4704.2142
$ws.Range("L137").Value = 40717.125
$ws.Range("M137").Value = -2154.2142
$ws.Range("N137").Value = -45817.125

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4502.6597
$ws.Range("I32").Value = 4502.6597
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4502.6597
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = -4215.6597

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 17506.562
$ws.Range("J45").Value = 4880.5557
$ws.Range("L45").Value = 4880.5557
$ws.Range("N45").Value = -5634.5557

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3561.7666
$ws.Range("I61").Value = 2197.182
$ws.Range("K61").Value = 2197.182
$ws.Range("M61").Value = -1985.182

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 316152.06
$ws.Range("I74").Value = 506556.72
$ws.Range("J74").Value = 16944.715
$ws.Range("K74").Value = 506556.72
$ws.Range("L74").Value = 16944.715
$ws.Range("M74").Value = -505682.72
$ws.Range("N74").Value = -18692.715

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 316152.06
$ws.Range("I77").Value = 506556.72
$ws.Range("J77").Value = 16944.715
$ws.Range("K77").Value = 2532783.6
$ws.Range("L77").Value = 84723.575
$ws.Range("M77").Value = -2528415.6
$ws.Range("N77").Value = -93459.575

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4371.6895
$ws.Range("I110").Value = 3722.7693
$ws.Range("K110").Value = 3722.7693
$ws.Range("M110").Value = -1677.7693

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3561.7666
$ws.Range("I136").Value = 2197.182
$ws.Range("K136").Value = 6591.545999999999
$ws.Range("M136").Value = -4041.545999999999

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15297145
$ws.Range("I105").Value = 835980.7
$ws.Range("J105").Value = 50003940
$ws.Range("K105").Value = 835980.7
$ws.Range("L105").Value = 50003940
$ws.Range("M105").Value = -834233.7
$ws.Range("N105").Value = -50007434

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2138016.5
$ws.Range("I107").Value = 2653577
$ws.Range("K107").Value = 2653577
$ws.Range("M107").Value = -2651657

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1722.8718
$ws.Range("I16").Value = 1659.4642
$ws.Range("K16").Value = 1659.4642
$ws.Range("M16").Value = -1372.4642

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3188.926
$ws.Range("I58").Value = 2274.5715
$ws.Range("J58").Value = 4173.615
$ws.Range("K58").Value = 2274.5715
$ws.Range("L58").Value = 4173.615
$ws.Range("M58").Value = -2071.5715
$ws.Range("N58").Value = -4579.615

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3086.75
$ws.Range("J105").Value = 3086.75
$ws.Range("L105").Value = 3086.75
$ws.Range("N105").Value = -6580.75

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2381402
$ws.Range("I107").Value = 2778219
$ws.Range("K107").Value = 2778219
$ws.Range("M107").Value = -2776299

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1722.8718
$ws.Range("I113").Value = 1659.4642
$ws.Range("K113").Value = 1659.4642
$ws.Range("M113").Value = 510.5358000000001

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14288526
$ws.Range("I132").Value = 20835722
$ws.Range("J132").Value = 3732.818
$ws.Range("K132").Value = 62507166
$ws.Range("L132").Value = 11198.454
$ws.Range("M132").Value = -62504636
$ws.Range("N132").Value = -16258.454

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4629.8887
$ws.Range("I134").Value = 4395.136
$ws.Range("K134").Value = 13185.408
$ws.Range("M134").Value = -10650.408

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3188.926
$ws.Range("I136").Value = 2274.5715
$ws.Range("J136").Value = 4173.615
$ws.Range("K136").Value = 6823.7145
$ws.Range("L136").Value = 12520.845
$ws.Range("M136").Value = -4273.7145
$ws.Range("N136").Value = -17620.845

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1937.3182
$ws.Range("J5").Value = 2243.889
$ws.Range("L5").Value = 6731.667
$ws.Range("N5").Value = -6955.667

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6522.857
$ws.Range("I113").Value = 495.33334
$ws.Range("K113").Value = 1486.00002
$ws.Range("M113").Value = 683.9999800000001

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1937.3182
$ws.Range("J135").Value = 2243.889
$ws.Range("L135").Value = 20195.001
$ws.Range("N135").Value = -25265.001

# Sheet CUL, row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 11757
$ws.Range("I141").Value = 11757
$ws.Range("K141").Value = 35271
$ws.Range("M141").Value = -30091

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3313.1843
$ws.Range("I113").Value = 3135.147
$ws.Range("K113").Value = 3135.147
$ws.Range("M113").Value = -965.1469999999999

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9971.4375
$ws.Range("I122").Value = 8954.6
$ws.Range("K122").Value = 26863.8
$ws.Range("M122").Value = -24413.8

# Sheet GSM, row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 102999
$ws.Range("J138").Value = 102999
$ws.Range("L138").Value = 102999
$ws.Range("N138").Value = -113279

# Sheet GSM, row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 99927.164
$ws.Range("J139").Value = 99927.164
$ws.Range("L139").Value = 99927.164
$ws.Range("N139").Value = -110207.164

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3961.1724
$ws.Range("I7").Value = 3313.5217
$ws.Range("K7").Value = 3313.5217
$ws.Range("M7").Value = -3201.5217

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2181.1538
$ws.Range("I61").Value = 2172.2
$ws.Range("K61").Value = 2172.2
$ws.Range("M61").Value = -1970.2

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2181.1538
$ws.Range("I113").Value = 2172.2
$ws.Range("K113").Value = 2172.2
$ws.Range("M113").Value = -2.199999999999818

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3961.1724
$ws.Range("I126").Value = 3313.5217
$ws.Range("K126").Value = 9940.5651
$ws.Range("M126").Value = -7470.5651

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3339.138
$ws.Range("I132").Value = 2663.48
$ws.Range("J132").Value = 7562
$ws.Range("K132").Value = 7990.440000000001
$ws.Range("L132").Value = 22686
$ws.Range("M132").Value = -5460.440000000001
$ws.Range("N132").Value = -27746

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3418.3777
$ws.Range("I136").Value = 3222.7097
$ws.Range("J136").Value = 3851.6428
$ws.Range("K136").Value = 9668.1291
$ws.Range("L136").Value = 11554.9284
$ws.Range("M136").Value = -7118.1291
$ws.Range("N136").Value = -16654.9284

# Sheet WVR, row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5901.421
$ws.Range("J96").Value = 5853.857
$ws.Range("L96").Value = 5853.857
$ws.Range("N96").Value = -8599.857

# Sheet WVR, row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 41667520
$ws.Range("I100").Value = 757.7857
$ws.Range("K100").Value = 1515.5714
$ws.Range("M100").Value = -974.5714

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 475.07144
$ws.Range("I107").Value = 470.4
$ws.Range("J107").Value = 486.75
$ws.Range("K107").Value = 1411.2
$ws.Range("L107").Value = 1460.25
$ws.Range("M107").Value = 508.8000000000002
$ws.Range("N107").Value = -5300.25

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 636.05
$ws.Range("J113").Value = 725.4286
$ws.Range("L113").Value = 2176.2858
$ws.Range("N113").Value = -6516.2858

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8774705
$ws.Range("I132").Value = 12348456
$ws.Range("J132").Value = 2769.9092
$ws.Range("K132").Value = 37045368
$ws.Range("L132").Value = 8309.7276
$ws.Range("M132").Value = -37042838
$ws.Range("N132").Value = -13369.7276

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21278578
$ws.Range("I136").Value = 23810404
$ws.Range("J136").Value = 11248.6
$ws.Range("K136").Value = 71431212
$ws.Range("L136").Value = 33745.8
$ws.Range("M136").Value = -71428662
$ws.Range("N136").Value = -38845.8
